# feat: add 2022-Q1 data
#
# 1. Duplicate the "总计" sheet (it already carries the header/index-column
#    style we need) to create a new "2022-Q1" sheet positioned right after
#    "2021-Q3", then overwrite its contents with the new fund-holdings data.
# 2. Insert a new row at the top of the data in "总计" for the 2022-Q1
#    summary, pushing the existing 2021-Q3 summary row down.
#
# NOTE: worksheet object references captured before a Copy()/order-changing
# call can go stale (end up pointing at the wrong sheet) - so sheets are
# re-fetched by name immediately after the Copy() below instead of reusing
# the pre-copy variable.

$wb = $excel.ActiveWorkbook

$wsQ3 = $wb.Worksheets.Item(1)
$origTotalName = $wb.Worksheets.Item(2).Name

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet as a copy of "总计" so it inherits
# the bordered/bold style already used for the header row + index column.
# ---------------------------------------------------------------------
$wb.Worksheets.Item($origTotalName).Copy($null, $wsQ3)

# The copy lands right after "2021-Q3" (position 2); re-fetch it fresh by
# position and rename it, then re-fetch the original "总计" sheet by name
# (it has shifted down to position 3).
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Item($origTotalName)

# Extend the styled header row (currently B1:D1) across to H1.
$wsNew.Range("D1").Copy()
$wsNew.Range("E1:H1").PasteSpecial(-4122)

# Extend the styled index column (currently A2) down through A4.
$wsNew.Range("A2").Copy()
$wsNew.Range("A3:A4").PasteSpecial(-4122)

# Headers
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Index column values (row below header = 0, 1, 2 ...)
$wsNew.Range("A2").Value = 0
$wsNew.Range("A3").Value = 1
$wsNew.Range("A4").Value = 2

# Data rows - text-like columns keep their values as text (quote-prefix
# forces text typing, ClearFormats then strips the quote-prefix style so
# the cell ends up unstyled, matching the source data's plain cells).
$textCells = @{
    "B2" = "006199";   "C2" = "长盛同锦研究精选混合";               "D2" = "1.73"; "E2" = "82.48"; "F2" = "3.51"; "G2" = "0.0607"
    "B3" = "001892";   "C3" = "长盛新兴成长主题灵活配置混合";       "D3" = "1.32"; "E3" = "82.10"; "F3" = "3.51"; "G3" = "0.0463"
    "B4" = "002085";   "C4" = "长盛互联网+主题灵活配置混合";       "D4" = "0.84"; "E4" = "83.97"; "F4" = "3.47"; "G4" = "0.0291"
}
foreach ($addr in $textCells.Keys) {
    $cell = $wsNew.Range($addr)
    $cell.Value = "'" + $textCells[$addr]
    $cell.ClearFormats()
}

# Numeric rank column
$wsNew.Range("H2").Value = 2
$wsNew.Range("H3").Value = 3
$wsNew.Range("H4").Value = 4

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q1 summary row above the existing 2021-Q3 row
# in "总计", then fix up the index numbers.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.14

$wsTotal.Range("A3").Value = 1
